# Applies the Renaissance Shenandoah-GC heap-8G stats refresh:
#   - the 4 summary cells (rows 1-4) are replaced with the new headline
#     values ("0M", "0M", "0M", "1026")
#   - the per-phase timing cells (rows 6,7,8,9,10,11,12) get refreshed
#     numbers
#   - the last three rows (44-46), which each held a whole tab-separated
#     line of per-run numbers crammed into a single run, are collapsed
#     down to just their first (summary) figure
$d = $word.ActiveDocument
$t = $d.Tables(1)

$t.Cell(1, 1).Range.Text  = "0M"
$t.Cell(2, 1).Range.Text  = "0M"
$t.Cell(3, 1).Range.Text  = "0M"
$t.Cell(4, 1).Range.Text  = "1026"

$t.Cell(6, 1).Range.Text  = "0.00064"
$t.Cell(7, 1).Range.Text  = "0.00020"
$t.Cell(8, 1).Range.Text  = "0.00006"
$t.Cell(9, 1).Range.Text  = "0.00026"
$t.Cell(10, 1).Range.Text = "0.00031"
$t.Cell(11, 1).Range.Text = "0.00042"
$t.Cell(12, 1).Range.Text = "0.19875"

$t.Cell(44, 1).Range.Text = "99.9"
$t.Cell(45, 1).Range.Text = "0.2"
$t.Cell(46, 1).Range.Text = "191"
